$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.432.52'
$ws.Cells.Item(2, 5).Value = '  +1.22%  '
$ws.Cells.Item(3, 4).Value = '3.151.03'
$ws.Cells.Item(3, 5).Value = '  +0.89%  '
$ws.Cells.Item(4, 5).Value = '  +0.12%  '
$ws.Cells.Item(5, 4).Value = "'527.76"
$ws.Cells.Item(5, 5).Value = '  -0.83%  '
$ws.Cells.Item(6, 4).Value = "'139.55"
$ws.Cells.Item(6, 5).Value = '  +1.01%  '
$ws.Cells.Item(7, 5).Value = '  +0.03%  '
$ws.Cells.Item(8, 4).Value = "'0.529"
$ws.Cells.Item(8, 5).Value = '  +10.94%  '
$ws.Cells.Item(9, 4).Value = "'7.26"
$ws.Cells.Item(9, 5).Value = '  -0.76%  '
$ws.Cells.Item(10, 4).Value = "'0.438"
$ws.Cells.Item(10, 5).Value = '  +6.01%  '
$ws.Cells.Item(11, 4).Value = "'0.111"
$ws.Cells.Item(11, 5).Value = '  +3.80%  '
$ws.Cells.Item(12, 5).Value = '  +2.21%  '
$ws.Cells.Item(13, 4).Value = '3.707.10'
$ws.Cells.Item(13, 5).Value = '  +1.34%  '
$ws.Cells.Item(14, 4).Value = "'25.58"
$ws.Cells.Item(14, 5).Value = '  -1.45%  '
$ws.Cells.Item(15, 4).Value = "'0.0000169"
$ws.Cells.Item(15, 5).Value = '  +3.04%  '
$ws.Cells.Item(16, 4).Value = '58.575.22'
$ws.Cells.Item(16, 5).Value = '  +1.27%  '
$ws.Cells.Item(17, 4).Value = '3.165.83'
$ws.Cells.Item(17, 5).Value = '  +1.14%  '
$ws.Cells.Item(18, 4).Value = "'6.20"
$ws.Cells.Item(18, 5).Value = '  +2.24%  '
$ws.Cells.Item(19, 4).Value = "'12.90"
$ws.Cells.Item(19, 5).Value = '  +1.63%  '
$ws.Cells.Item(20, 4).Value = "'8.09"
$ws.Cells.Item(20, 5).Value = '  +0.18%  '
$ws.Cells.Item(21, 4).Value = "'372.69"
$ws.Cells.Item(21, 5).Value = '  +1.65%  '
$ws.Cells.Item(22, 4).Value = "'0.998"
$ws.Cells.Item(22, 5).Value = '  -0.05%  '
$ws.Cells.Item(23, 4).Value = "'0.526"
$ws.Cells.Item(23, 5).Value = '  +4.06%  '
$ws.Cells.Item(24, 4).Value = "'69.57"
$ws.Cells.Item(24, 5).Value = '  +0.65%  '
$ws.Cells.Item(25, 5).Value = '  +0.22%  '
$ws.Cells.Item(26, 4).Value = "'1.00"
$ws.Cells.Item(26, 5).Value = '  +0.06%  '
$ws.Cells.Item(27, 4).Value = "'8.40"
$ws.Cells.Item(27, 5).Value = '  +15.13%  '
$ws.Cells.Item(28, 4).Value = '0.0₃0848'
$ws.Cells.Item(28, 5).Value = '  -1.77%  '
$ws.Cells.Item(29, 4).Value = "'22.34"
$ws.Cells.Item(29, 5).Value = '  +4.45%  '
$ws.Cells.Item(30, 4).Value = "'1.87"
$ws.Cells.Item(30, 5).Value = '  +0.32%  '
$ws.Cells.Item(31, 4).Value = "'5.97"
$ws.Cells.Item(31, 5).Value = '  -1.54%  '
$ws.Cells.Item(32, 4).Value = "'5.11"
$ws.Cells.Item(32, 5).Value = '  -0.44%  '
$ws.Cells.Item(33, 4).Value = "'1.14"
$ws.Cells.Item(33, 5).Value = '  -0.07%  '
$ws.Cells.Item(34, 4).Value = "'6.28"
$ws.Cells.Item(34, 5).Value = '  +3.59%  '
$ws.Cells.Item(35, 4).Value = "'156.93"
$ws.Cells.Item(35, 5).Value = '  -1.67%  '
$ws.Cells.Item(36, 4).Value = "'1.33"
$ws.Cells.Item(36, 5).Value = '  +2.56%  '
$ws.Cells.Item(37, 4).Value = '2.692.13'
$ws.Cells.Item(37, 5).Value = '  +6.31%  '
$ws.Cells.Item(38, 4).Value = "'24.89"
$ws.Cells.Item(38, 5).Value = '  -2.08%  '
$ws.Cells.Item(39, 4).Value = "'1.68"
$ws.Cells.Item(39, 5).Value = '  +0.20%  '
$ws.Cells.Item(40, 4).Value = "'0.0688"
$ws.Cells.Item(40, 5).Value = '  +2.67%  '
$ws.Cells.Item(41, 4).Value = "'4.25"
$ws.Cells.Item(41, 5).Value = '  +5.50%  '
$ws.Cells.Item(42, 2).Value = 'Mantle'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(42, 4).Value = "'0.721"
$ws.Cells.Item(42, 5).Value = '  +3.20%  '
$ws.Cells.Item(43, 2).Value = 'VeChain'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(43, 4).Value = "'0.0291"
$ws.Cells.Item(43, 5).Value = '  +8.05%  '
$ws.Cells.Item(44, 4).Value = "'38.99"
$ws.Cells.Item(44, 5).Value = '  +3.18%  '
$ws.Cells.Item(45, 5).Value = '  +0.16%  '
$ws.Cells.Item(46, 4).Value = '3.200.92'
$ws.Cells.Item(46, 5).Value = '  +1.13%  '
$ws.Cells.Item(47, 4).Value = "'0.102"
$ws.Cells.Item(47, 5).Value = '  +11.46%  '
$ws.Cells.Item(48, 4).Value = "'6.18"
$ws.Cells.Item(48, 5).Value = '  +1.27%  '
$ws.Cells.Item(49, 4).Value = "'0.976"
$ws.Cells.Item(49, 5).Value = '  -0.47%  '
$ws.Cells.Item(50, 4).Value = "'19.96"
$ws.Cells.Item(50, 5).Value = '  +1.25%  '
$ws.Cells.Item(51, 4).Value = "'0.744"
$ws.Cells.Item(51, 5).Value = '  +0.72%  '
